$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 19, shifting existing rows 19-24 down to 20-25
$ws.Rows.Item(19).Insert()

# Fill the new row 19 with data
$ws.Cells.Item(19, 1).Value = 7
$ws.Cells.Item(19, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(19, 3).Value = "Ñuble"
$ws.Cells.Item(19, 4).Value = 44627
$ws.Cells.Item(19, 5).Value = 16
$ws.Cells.Item(19, 6).Value = 100112001
$ws.Cells.Item(19, 7).Value = "Berenjena"
$ws.Cells.Item(19, 8).Value = "Sin especificar"
$ws.Cells.Item(19, 9).Value = "Primera"
$ws.Cells.Item(19, 10).Value = 60
$ws.Cells.Item(19, 11).Value = 9000
$ws.Cells.Item(19, 12).Value = 9500
$ws.Cells.Item(19, 13).Value = 9250
$ws.Cells.Item(19, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(19, 15).Value = "Región Metropolitana"
$ws.Cells.Item(19, 16).Value = 154
$ws.Cells.Item(19, 17).Value = 60
$ws.Cells.Item(19, 18).Value = "Hortaliza"
